# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ)
# across the ALC / ARM / BSM / CRP / CUL / LTW / WVR leve-profit tables.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values,
        [string[]]$ClearCols = @()
    )
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
    foreach ($col in $ClearCols) {
        $ws.Range("$col$Row").ClearContents()
    }
}

# ---------------- ALC (sheet1) ----------------
Set-Row -SheetName "ALC" -Row 33 -Values @{
    H = 420.84
    I = 457.77274
    J = 150
    K = 457.77274
    L = 150
    M = -228.77274
    N = -608
}

Set-Row -SheetName "ALC" -Row 74 -Values @{
    H = 3004
    I = 0
    K = 0
} -ClearCols @("M")

Set-Row -SheetName "ALC" -Row 77 -Values @{
    H = 3004
    I = 0
    K = 0
} -ClearCols @("M")

Set-Row -SheetName "ALC" -Row 86 -Values @{
    H = 39299500
    I = 47842236
    K = 47842236
    M = -47841113
}

Set-Row -SheetName "ALC" -Row 89 -Values @{
    H = 39299500
    I = 47842236
    K = 239211180
    M = -239205564
}

Set-Row -SheetName "ALC" -Row 111 -Values @{
    H = 2456.5557
    I = 1304.8334
    J = 4760
    K = 3914.5002
    L = 14280
    M = -847.5001999999999
    N = -20414
}

Set-Row -SheetName "ALC" -Row 115 -Values @{
    H = 2623.125
    I = 2246.25
    J = 3000
    K = 6738.75
    L = 9000
    M = -5171.75
    N = -12134
}

Set-Row -SheetName "ALC" -Row 116 -Values @{
    H = 2482.0908
    I = 2280.5557
    J = 3389
    K = 2280.5557
    L = 3389
    M = 1161.4443
    N = -10273
}

# ---------------- ARM (sheet2) ----------------
Set-Row -SheetName "ARM" -Row 132 -Values @{
    H = 3761.8704
    I = 3334.475
    K = 10003.425
    M = -7473.424999999999
}

# ---------------- BSM (sheet3) ----------------
Set-Row -SheetName "BSM" -Row 22 -Values @{
    H = 590.7143
    I = 590.7143
    K = 590.7143
    M = -417.7143
}

Set-Row -SheetName "BSM" -Row 64 -Values @{
    H = 480.25
    I = 419.5
    J = 601.75
    K = 419.5
    L = 601.75
    M = -194.5
    N = -1051.75
}

Set-Row -SheetName "BSM" -Row 67 -Values @{
    H = 480.25
    I = 419.5
    J = 601.75
    K = 419.5
    L = 601.75
    M = 360.5
    N = -2161.75
}

# ---------------- CRP (sheet4) ----------------
Set-Row -SheetName "CRP" -Row 58 -Values @{
    H = 3100
    I = 2900
    J = 3500
    K = 2900
    L = 3500
    M = -2697
    N = -3906
}

Set-Row -SheetName "CRP" -Row 105 -Values @{
    H = 839.3333
    I = 806.1539
    K = 806.1539
    M = 940.8461
}

Set-Row -SheetName "CRP" -Row 136 -Values @{
    H = 3100
    I = 2900
    J = 3500
    K = 8700
    L = 10500
    M = -6150
    N = -15600
}

Set-Row -SheetName "CRP" -Row 141 -Values @{
    H = 118512.38
    J = 119724.21
    L = 119724.21
    N = -130084.21
}

# ---------------- CUL (sheet5) ----------------
Set-Row -SheetName "CUL" -Row 17 -Values @{
    H = 1214.9
    I = 807.1429000000001
    J = 2166.3333
    K = 2421.4287
    L = 6498.999899999999
    M = -2252.4287
    N = -6836.999899999999
}

Set-Row -SheetName "CUL" -Row 34 -Values @{
    H = 8197293
    J = 8621286
    L = 25863858
    N = -25864026
}

Set-Row -SheetName "CUL" -Row 39 -Values @{
    H = 2598.4285
    J = 2598.4285
    L = 7795.2855
    N = -8383.2855
}

Set-Row -SheetName "CUL" -Row 55 -Values @{
    H = 1291.2916
    J = 1291.2916
    L = 3873.8748
    N = -4227.8748
}

Set-Row -SheetName "CUL" -Row 131 -Values @{
    H = 3058.3774
    I = 409.9375
    J = 4203.6484
    K = 1229.8125
    L = 12610.9452
    M = 3810.1875
    N = -22690.9452
}

Set-Row -SheetName "CUL" -Row 137 -Values @{
    H = 30353.85
    J = 53698.35
    L = 161095.05
    N = -171295.05
}

Set-Row -SheetName "CUL" -Row 141 -Values @{
    H = 8679.947
    I = 7509.9165
    J = 10685.714
    K = 22529.7495
    L = 32057.142
    M = -17349.7495
    N = -42417.142
}

# ---------------- LTW (sheet7) ----------------
Set-Row -SheetName "LTW" -Row 61 -Values @{
    H = 5209
    I = 4627
    J = 5500
    K = 4627
    L = 5500
    M = -4425
    N = -5904
}

Set-Row -SheetName "LTW" -Row 82 -Values @{
    H = 2080.111
    I = 1700
    J = 2226.3076
    K = 1700
    L = 2226.3076
    M = -1339
    N = -2948.3076
}

Set-Row -SheetName "LTW" -Row 85 -Values @{
    H = 2080.111
    I = 1700
    J = 2226.3076
    K = 1700
    L = 2226.3076
    M = -452
    N = -4722.3076
}

Set-Row -SheetName "LTW" -Row 113 -Values @{
    H = 5209
    I = 4627
    J = 5500
    K = 4627
    L = 5500
    M = -2457
    N = -9840
}

Set-Row -SheetName "LTW" -Row 141 -Values @{
    H = 70000
    J = 70000
    L = 70000
    N = -80360
}

# ---------------- WVR (sheet8) ----------------
Set-Row -SheetName "WVR" -Row 113 -Values @{
    H = 1803.0834
    I = 1975
    J = 943.5
    K = 5925
    L = 2830.5
    M = -3755
    N = -7170.5
}

Set-Row -SheetName "WVR" -Row 136 -Values @{
    H = 3294.4443
    I = 2727.3438
    J = 4690.385
    K = 8182.0314
    L = 14071.155
    M = -5632.0314
    N = -19171.155
}

Write-Host "Anima_Profits scheduled update applied"
